# GSC export refresh: slide the 90-day window forward by one day.
# Drop the oldest date row (2025-10-24), shift every remaining row up by
# one, and append the new day (2026-01-22) with its HTTPS URL count.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$lastRow = 91

# Shift rows 3..91 up into rows 2..90 (A:C) in one bulk copy so the
# existing shared-string cells (dates) are relocated as-is instead of
# being re-typed (re-typing a "yyyy-MM-dd" looking string makes Excel
# "helpfully" reinterpret it as a real date and reformat the cell).
$srcRange = $ws.Range("A3:C" + $lastRow)
$dstRange = $ws.Range("A2:C" + ($lastRow - 1))
$srcRange.Copy($dstRange)

# Fill in the newly-opened last row with the new day.
$newDateCell = $ws.Cells.Item($lastRow, 1)

# Write the date through a formula + Paste-Special-Values round trip so
# it lands as plain text (matches the rest of column A) instead of
# letting the "=" assignment path auto-convert the literal into an
# Excel date serial/date-formatted cell.
$newDateCell.Formula = '="2026-01-22"'
$newDateCell.Copy()
$newDateCell.PasteSpecial(-4163)
$excel.CutCopyMode = 0

$ws.Cells.Item($lastRow, 2).Value = 0.0
$ws.Cells.Item($lastRow, 3).Value = 24.0
